$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update CuahsiUsed (column G) from "yes" to "no" for rows 9, 10, 18 and 19
$ws.Cells.Item(9, 7).Value = "no"
$ws.Cells.Item(10, 7).Value = "no"
$ws.Cells.Item(18, 7).Value = "no"
$ws.Cells.Item(19, 7).Value = "no"

# Update the active selection to reflect where the editing ended up (G19)
$ws.Range("G19").Select()
